{"js": "// Afrikaans translation of the \"Child Safety interview information & consent\n// form\" document. Replaces the English source text with the Afrikaans\n// translation, run by run, using Word's search-and-replace (the runs keep\n// their original formatting because insertText(..., Replace) only touches\n// the text content of the matched range).\n\nconst replacements = [\n  {\n    find:\n      \"Appendix 19: SWIFT Child Safety Module Interview: Information Sheet and Consent Form\",\n    replace:\n      \"Bylae 19: SWIFT Kindersveiligheidsmodule Onderhoud: Inligtingsblad en Toestemmingsvorm\",\n  },\n  {\n    find:\n      \"Your experience with the ParentText programme is vital to our study. \" +\n      \"We would love to hear about your experience of the \\u2018Keeping My Child Safe\\u2019 module. \" +\n      \"This interview is part of a study carried out by researchers from the Universities of Cape Town \" +\n      \"in South Africa and the University of Oxford in the United Kingdom. \",\n    replace:\n      \"Jou ervaring met die ParentText-program is van groot belang vir ons studie. \" +\n      \"Ons sal graag wil hoor oor jou ervaring met die \\u201cHou My Kind Veilig\\u201d module. \" +\n      \"Hierdie onderhoud is deel van 'n studie wat uitgevoer word deur navorsers van die \" +\n      \"Universiteit van Kaapstad in Suid-Afrika en die Universiteit van Oxford in die Verenigde Koninkryk. \",\n  },\n  {\n    find:\n      \"Before you decide if you\\u2019d like to be interviewed, it\\u2019s important for you to know why \" +\n      \"we\\u2019re doing this research and what participating in it would involve. All the information \" +\n      \"you might need is explained below but if you have any questions about your participation or our \" +\n      \"study, please email the study team at swift@globalparenting.org or message us on WhatsApp at \" +\n      \"+27 XX XXX XXXX. We\\u2019re here to help you!\",\n    replace:\n      \"Voordat jy besluit of jy aan die onderhoud wil deelneem, is dit belangrik om te weet waarom ons \" +\n      \"hierdie navorsing doen en wat deelname behels. Alle inligting wat jy mag benodig, word hieronder \" +\n      \"verduidelik, maar as jy enige vrae het oor jou deelname of ons studie, kan jy die studiespan per \" +\n      \"e-pos kontak by swift@globalparenting.org of 'n boodskap stuur op WhatsApp na +27 XX XXX XXXX. \" +\n      \"Ons is hier om jou te help!\",\n  },\n  {\n    find: \" What will my interview look like and what is expected of me?\",\n    replace: \" Hoe sal my onderhoud lyk en wat word van my verwag?\",\n  },\n  {\n    find:\n      \"We would like to have a telephonic conversation with you which will last a maximum of 45 minutes. \" +\n      \"One of our researchers will call you to speak to you at a time that is convenient for you. \" +\n      \"There are no right or wrong answers, we just want to hear your experience and opinion of the chatbot. \" +\n      \"Please make sure that when we call, that you only let the interview start when you are in a private \" +\n      \"space where you feel comfortable to talk without being overheard or interrupted. If while you are \" +\n      \"being interviewed, you are interrupted, please ask the researcher to pause until you feel safe to \" +\n      \"continue talking.\",\n    replace:\n      \"Ons wil graag 'n telefoniese gesprek met jou voer wat 'n maksimum van 45 minute sal duur. \" +\n      \"One of our researchers will call you to speak to you at a time that is convenient for you. \" +\n      \"Daar is geen regte of verkeerde antwoorde nie; ons wil net jou ervaring en mening oor die geselsbot \" +\n      \"hoor. Please make sure that when we call, that you only let the interview start when you are in a \" +\n      \"private space where you feel comfortable to talk without being overheard or interrupted. If while \" +\n      \"you are being interviewed, you are interrupted, please ask the researcher to pause until you feel \" +\n      \"safe to continue talking.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + find.substring(0, 60));\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Afrikaans translation of the \"Child Safety interview information & consent\n# form\" document. Replaces the English source text with the Afrikaans\n# translation using Find/Execute to locate each passage and then assigning\n# the matched Range's .Text directly (NOT Find.Replacement.Text) so that\n# Word's smart-quotes autocorrect does not mangle straight apostrophes\n# (e.g. \"'n\") into curly opening quotes.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $found = $r.Find.Execute($findText, $false, $true)\n    if (-not $found) {\n        throw \"Could not find expected text: $findText\"\n    }\n    $r.Text = $replaceText\n}\n\nReplace-Text `\n    'Appendix 19: SWIFT Child Safety Module Interview: Information Sheet and Consent Form' `\n    'Bylae 19: SWIFT Kindersveiligheidsmodule Onderhoud: Inligtingsblad en Toestemmingsvorm'\n\nReplace-Text `\n    'Your experience with the ParentText programme is vital to our study. We would love to hear about your experience of the \u2018Keeping My Child Safe\u2019 module. This interview is part of a study carried out by researchers from the Universities of Cape Town in South Africa and the University of Oxford in the United Kingdom. ' `\n    'Jou ervaring met die ParentText-program is van groot belang vir ons studie. Ons sal graag wil hoor oor jou ervaring met die \u201cHou My Kind Veilig\u201d module. Hierdie onderhoud is deel van ''n studie wat uitgevoer word deur navorsers van die Universiteit van Kaapstad in Suid-Afrika en die Universiteit van Oxford in die Verenigde Koninkryk. '\n\nReplace-Text `\n    'Before you decide if you\u2019d like to be interviewed, it\u2019s important for you to know why we\u2019re doing this research and what participating in it would involve. All the information you might need is explained below but if you have any questions about your participation or our study, please email the study team at swift@globalparenting.org or message us on WhatsApp at +27 XX XXX XXXX. We\u2019re here to help you!' `\n    'Voordat jy besluit of jy aan die onderhoud wil deelneem, is dit belangrik om te weet waarom ons hierdie navorsing doen en wat deelname behels. Alle inligting wat jy mag benodig, word hieronder verduidelik, maar as jy enige vrae het oor jou deelname of ons studie, kan jy die studiespan per e-pos kontak by swift@globalparenting.org of ''n boodskap stuur op WhatsApp na +27 XX XXX XXXX. Ons is hier om jou te help!'\n\nReplace-Text `\n    ' What will my interview look like and what is expected of me?' `\n    ' Hoe sal my onderhoud lyk en wat word van my verwag?'\n\nReplace-Text `\n    'We would like to have a telephonic conversation with you which will last a maximum of 45 minutes. One of our researchers will call you to speak to you at a time that is convenient for you. There are no right or wrong answers, we just want to hear your experience and opinion of the chatbot. Please make sure that when we call, that you only let the interview start when you are in a private space where you feel comfortable to talk without being overheard or interrupted. If while you are being interviewed, you are interrupted, please ask the researcher to pause until you feel safe to continue talking.' `\n    'Ons wil graag ''n telefoniese gesprek met jou voer wat ''n maksimum van 45 minute sal duur. One of our researchers will call you to speak to you at a time that is convenient for you. Daar is geen regte of verkeerde antwoorde nie; ons wil net jou ervaring en mening oor die geselsbot hoor. Please make sure that when we call, that you only let the interview start when you are in a private space where you feel comfortable to talk without being overheard or interrupted. If while you are being interviewed, you are interrupted, please ask the researcher to pause until you feel safe to continue talking.'\n"}
